$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A56").Value = "A reel showcasing Delhi University’s famous Paneer Singapore Chowmein, priced at 180/-, located at Nams Momos, Kamla Nagar, Delhi. The content focuses on food and local Delhi culture, with high engagement from an Indian audience."
$ws.Range("B56").Value = "Food enthusiast and local Delhi explorer, likely a student or young adult, creating relatable food content for a Delhi-based audience."
$ws.Range("C56").Value = "Casual and appreciative, focusing on the food and local vibe."
$ws.Range("D56").Value = "Medium"
$ws.Range("E56").Value = "both"
$ws.Range("F56").Value = 7
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = "https://www.instagram.com/p/DL4bFQiyOhX/"
$ws.Range("I56").Value = "https://www.instagram.com/yumm_hmm__"
$ws.Range("J56").Value = "The creator focuses on Delhi-specific food content, indicating a strong Indian audience, which aligns with Knytt's target market. The reel has positive engagement (emojis like 😍, 🔥, ❤️), suggesting an active and interested audience. While the creator may be interested in joining a platform like Knytt to connect with foodies, the likelihood is medium due to the niche focus on food rather than personal interaction. A comment can engage the audience, while a personalized DM can directly pitch to the creator."
$ws.Range("K56").Value = "both"
$ws.Range("L56").Value = "Hey @yumm_hmm__, loved your reel on Paneer Singapore Chowmein! 😍 We’re building a cool app called Knytt where foodies like you can connect with fans through texting and video calls. Would love to have you on board to share your Delhi food adventures! Check us out at knytt.in or drop me a message if you’re curious. ✨"
$ws.Range("M56").Value = "This looks insanely delicious! 😍 Kamla Nagar always has the best eats. If you’re a foodie, you’d love connecting with more peeps on @knyttofficial! 🍝🔥"
